# Update countries & provincias Spain
# Applies the country-ranking reshuffle + refreshed case counts described by the
# commit "Update countries & provincias Spain".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row([int]$Row, [string]$Country, [double]$B, [double]$C, [double]$D, [double]$E, [double]$F, [double]$G, [double]$H) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
}

# Row 4: Estados Unidos - refreshed totals
Set-Row 4 "Estados Unidos" 699850 22280 59142 603467 13466 2624 37241

# Row 8: Alemania - refreshed totals
Set-Row 8 "Alemania" 140886 3188 83114 53446 5013 274 4326

# Rows 21/22: India overtakes Irlanda
Set-Row 21 "India" 14352 922 2041 11825 0 38 486
Set-Row 22 "Irlanda" 13980 709 77 13373 156 44 530

# Rows 27/28: Japon overtakes Chile
Set-Row 27 "Japon" 9787 556 935 8662 193 0 190
Set-Row 28 "Chile" 9252 445 3621 5515 384 11 116

# Row 52: Colombia - refreshed totals
Set-Row 52 "Colombia" 3439 206 634 2652 98 9 153

# Rows 84/85: Tunez overtakes Bulgaria
Set-Row 84 "Tunez" 864 42 43 784 89 0 37
Set-Row 85 "Bulgaria" 846 46 141 664 37 3 41

# Row 96: Burkina Faso - refreshed totals
Set-Row 96 "Burkina Faso" 557 11 294 228 0 3 35

# Rows 99-103: Nigeria moves up ahead of Kirguistan; Guinea, Bolivia, Honduras shift down
Set-Row 99 "Nigeria" 493 51 159 317 2 4 17
Set-Row 100 "Kirguistan" 489 23 114 370 5 0 5
Set-Row 101 "Guinea" 477 39 59 415 0 2 3
Set-Row 102 "Bolivia" 465 24 26 408 3 2 31
Set-Row 103 "Honduras" 442 16 10 391 10 6 41

# Rows 154/155: Uganda overtakes Polinesia Francesa
Set-Row 154 "Uganda" 56 1 20 36 0 0 0
Set-Row 155 "Polinesia Francesa" 55 0 0 55 1 0 0
